$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24
$ws.Range("F2").Value = 13
$ws.Range("H2").Value = 13

# Row 14 updates
$ws.Range("F14").Value = 17
$ws.Range("H14").Value = 17
